$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1815068493150685
$ws.Range("C2").Value = 0.571917808219178
$ws.Range("J2").Value = 0.00684931506849315
$ws.Range("P2").Value = 0.1095890410958904
$ws.Range("S2").Value = 0.1301369863013699
$ws.Range("B3").Value = 0.01162790697674419
$ws.Range("C3").Value = 0.03488372093023256
$ws.Range("J3").Value = 0.02906976744186046
$ws.Range("P3").Value = 0.7209302325581395
$ws.Range("S3").Value = 0.2034883720930233
$ws.Range("J4").Value = 0.09302325581395349
$ws.Range("P4").Value = 0.7209302325581395
$ws.Range("S4").Value = 0.186046511627907
$ws.Range("B6").Value = 0.06796116504854369
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.05825242718446602
$ws.Range("J6").Value = 0.2621359223300971
$ws.Range("O6").Value = 0.01456310679611651
$ws.Range("Q6").Value = 0.1213592233009709
$ws.Range("R6").Value = 0.07766990291262135
$ws.Range("S6").Value = 0.3883495145631068
$ws.Range("B7").Value = 0.1045454545454545
$ws.Range("D7").Value = 0.02272727272727273
$ws.Range("F7").Value = 0.02272727272727273
$ws.Range("J7").Value = 0.1590909090909091
$ws.Range("O7").Value = 0.00909090909090909
$ws.Range("Q7").Value = 0.2045454545454546
$ws.Range("R7").Value = 0.05909090909090909
$ws.Range("S7").Value = 0.4181818181818182
$ws.Range("B8").Value = 0.103950103950104
$ws.Range("D8").Value = 0.02494802494802495
$ws.Range("F8").Value = 0.06860706860706861
$ws.Range("J8").Value = 0.1081081081081081
$ws.Range("O8").Value = 0.0103950103950104
$ws.Range("Q8").Value = 0.1767151767151767
$ws.Range("R8").Value = 0.0893970893970894
$ws.Range("S8").Value = 0.4178794178794179
$ws.Range("B9").Value = 0.1367521367521368
$ws.Range("D9").Value = 0.01282051282051282
$ws.Range("F9").Value = 0.03846153846153846
$ws.Range("J9").Value = 0.09401709401709402
$ws.Range("O9").Value = 0.02991452991452992
$ws.Range("Q9").Value = 0.1752136752136752
$ws.Range("R9").Value = 0.05982905982905983
$ws.Range("S9").Value = 0.452991452991453
$ws.Range("B10").Value = 0.1129186602870813
$ws.Range("D10").Value = 0.02105263157894737
$ws.Range("F10").Value = 0.0861244019138756
$ws.Range("J10").Value = 0.1196172248803828
$ws.Range("O10").Value = 0.01913875598086124
$ws.Range("Q10").Value = 0.1894736842105263
$ws.Range("R10").Value = 0.0660287081339713
$ws.Range("S10").Value = 0.3856459330143541
$ws.Range("G11").Value = 0.1715210355987055
$ws.Range("J11").Value = 0.07443365695792881
$ws.Range("K11").Value = 0.226537216828479
$ws.Range("L11").Value = 0.517799352750809
$ws.Range("S11").Value = 0.009708737864077669
$ws.Range("G12").Value = 0.7988165680473372
$ws.Range("J12").Value = 0.1479289940828402
$ws.Range("L12").Value = 0.04733727810650887
$ws.Range("S12").Value = 0.005917159763313609
$ws.Range("F13").Value = 0.01754385964912281
$ws.Range("G13").Value = 0.6842105263157895
$ws.Range("J13").Value = 0.2280701754385965
$ws.Range("S13").Value = 0.07017543859649122
$ws.Range("F15").Value = 0.01923076923076923
$ws.Range("H15").Value = 0.1730769230769231
$ws.Range("I15").Value = 0.09615384615384616
$ws.Range("J15").Value = 0.2980769230769231
$ws.Range("K15").Value = 0.07211538461538461
$ws.Range("M15").Value = 0.009615384615384616
$ws.Range("O15").Value = 0.04326923076923077
$ws.Range("S15").Value = 0.2884615384615384
$ws.Range("F16").Value = 0.02717391304347826
$ws.Range("H16").Value = 0.1739130434782609
$ws.Range("I16").Value = 0.09782608695652174
$ws.Range("J16").Value = 0.3478260869565217
$ws.Range("K16").Value = 0.1304347826086956
$ws.Range("M16").Value = 0.01630434782608696
$ws.Range("O16").Value = 0.02717391304347826
$ws.Range("S16").Value = 0.1793478260869565
$ws.Range("F17").Value = 0.02061855670103093
$ws.Range("H17").Value = 0.2036082474226804
$ws.Range("I17").Value = 0.1443298969072165
$ws.Range("J17").Value = 0.3402061855670103
$ws.Range("K17").Value = 0.07731958762886598
$ws.Range("M17").Value = 0.02577319587628866
$ws.Range("O17").Value = 0.05412371134020619
$ws.Range("S17").Value = 0.134020618556701
$ws.Range("F18").Value = 0.03225806451612903
$ws.Range("H18").Value = 0.1935483870967742
$ws.Range("I18").Value = 0.1225806451612903
$ws.Range("J18").Value = 0.2967741935483871
$ws.Range("K18").Value = 0.1096774193548387
$ws.Range("M18").Value = 0.02580645161290323
$ws.Range("N18").Value = 0.006451612903225806
$ws.Range("O18").Value = 0.07096774193548387
$ws.Range("S18").Value = 0.1419354838709677
$ws.Range("F19").Value = 0.01145912910618793
$ws.Range("H19").Value = 0.2337662337662338
$ws.Range("I19").Value = 0.0932009167303285
$ws.Range("J19").Value = 0.3086325439266616
$ws.Range("K19").Value = 0.1161191749427044
$ws.Range("M19").Value = 0.02902979373567609
$ws.Range("N19").Value = 0.0007639419404125286
$ws.Range("O19").Value = 0.0718105423987777
$ws.Range("S19").Value = 0.1352177234530176
